# Update "想去人数" (column F) values on each sheet to match the new
# snapshot captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 2332
$ws1.Cells.Item(3, 6).Value = 493
$ws1.Cells.Item(5, 6).Value = 334
$ws1.Cells.Item(6, 6).Value = 334
$ws1.Cells.Item(7, 6).Value = 544
$ws1.Cells.Item(9, 6).Value = 760
$ws1.Cells.Item(11, 6).Value = 774
$ws1.Cells.Item(14, 6).Value = 390
$ws1.Cells.Item(16, 6).Value = 1020
$ws1.Cells.Item(17, 6).Value = 19965
$ws1.Cells.Item(18, 6).Value = 630
$ws1.Cells.Item(19, 6).Value = 66
$ws1.Cells.Item(20, 6).Value = 225
$ws1.Cells.Item(21, 6).Value = 282
$ws1.Cells.Item(23, 6).Value = 144
$ws1.Cells.Item(25, 6).Value = 8
$ws1.Cells.Item(26, 6).Value = 209
$ws1.Cells.Item(28, 6).Value = 330
$ws1.Cells.Item(29, 6).Value = 137

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 178
$ws2.Cells.Item(5, 6).Value = 83
$ws2.Cells.Item(6, 6).Value = 203
$ws2.Cells.Item(7, 6).Value = 223
$ws2.Cells.Item(8, 6).Value = 3392
$ws2.Cells.Item(10, 6).Value = 88
$ws2.Cells.Item(16, 6).Value = 3104

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 94
$ws3.Cells.Item(4, 6).Value = 571

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 94
$ws4.Cells.Item(5, 6).Value = 2332
$ws4.Cells.Item(6, 6).Value = 571
$ws4.Cells.Item(7, 6).Value = 493
$ws4.Cells.Item(9, 6).Value = 334
$ws4.Cells.Item(10, 6).Value = 334
$ws4.Cells.Item(11, 6).Value = 544
$ws4.Cells.Item(12, 6).Value = 178
$ws4.Cells.Item(15, 6).Value = 83
$ws4.Cells.Item(16, 6).Value = 203
$ws4.Cells.Item(18, 6).Value = 760
$ws4.Cells.Item(20, 6).Value = 774
$ws4.Cells.Item(23, 6).Value = 390
$ws4.Cells.Item(25, 6).Value = 1020
$ws4.Cells.Item(26, 6).Value = 19968
$ws4.Cells.Item(27, 6).Value = 223
$ws4.Cells.Item(28, 6).Value = 3392
$ws4.Cells.Item(30, 6).Value = 88
$ws4.Cells.Item(32, 6).Value = 630
$ws4.Cells.Item(33, 6).Value = 66
$ws4.Cells.Item(34, 6).Value = 225
$ws4.Cells.Item(37, 6).Value = 282
$ws4.Cells.Item(39, 6).Value = 144
$ws4.Cells.Item(41, 6).Value = 8
$ws4.Cells.Item(44, 6).Value = 209
$ws4.Cells.Item(46, 6).Value = 330
$ws4.Cells.Item(47, 6).Value = 137
$ws4.Cells.Item(48, 6).Value = 3105
